{"js": "// Replace the multiplication equations in the document's table cells with\n// their new values. Each old equation string is unique in the document, so a\n// simple exact-text search/replace per pair is reliable.\nconst replacements = [\n  [\"34\u00d765=2210\", \"99\u00d761=6039\"],\n  [\"26\u00d746=1196\", \"72\u00d779=5688\"],\n  [\"26\u00d754=1404\", \"63\u00d713=819\"],\n  [\"74\u00d794=6956\", \"65\u00d776=4940\"],\n  [\"79\u00d773=5767\", \"49\u00d764=3136\"],\n  [\"60\u00d726=1560\", \"58\u00d764=3712\"],\n  [\"35\u00d796=3360\", \"20\u00d769=1380\"],\n  [\"11\u00d790=990\", \"75\u00d779=5925\"],\n  [\"44\u00d777=3388\", \"91\u00d796=8736\"],\n  [\"47\u00d774=3478\", \"99\u00d773=7227\"],\n  [\"65\u00d752=3380\", \"24\u00d760=1440\"],\n  [\"46\u00d729=1334\", \"88\u00d755=4840\"],\n  [\"90\u00d733=2970\", \"60\u00d746=2760\"],\n  [\"46\u00d721=966\", \"90\u00d775=6750\"],\n  [\"29\u00d720=580\", \"54\u00d733=1782\"],\n  [\"56\u00d756=3136\", \"60\u00d788=5280\"],\n  [\"77\u00d761=4697\", \"12\u00d783=996\"],\n  [\"63\u00d737=2331\", \"48\u00d756=2688\"],\n  [\"46\u00d727=1242\", \"71\u00d760=4260\"],\n  [\"89\u00d714=1246\", \"42\u00d789=3738\"],\n  [\"71\u00d787=6177\", \"62\u00d735=2170\"],\n  [\"74\u00d740=2960\", \"99\u00d751=5049\"],\n  [\"49\u00d751=2499\", \"32\u00d729=928\"],\n  [\"45\u00d745=2025\", \"13\u00d794=1222\"],\n  [\"62\u00d750=3100\", \"52\u00d711=572\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication equations in the document's table cells with\n# their new values. Each old equation string is unique in the document, so a\n# simple Find/Replace per pair (scoped to the whole document content) is\n# reliable.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"34\u00d765=2210\", \"99\u00d761=6039\"),\n    @(\"26\u00d746=1196\", \"72\u00d779=5688\"),\n    @(\"26\u00d754=1404\", \"63\u00d713=819\"),\n    @(\"74\u00d794=6956\", \"65\u00d776=4940\"),\n    @(\"79\u00d773=5767\", \"49\u00d764=3136\"),\n    @(\"60\u00d726=1560\", \"58\u00d764=3712\"),\n    @(\"35\u00d796=3360\", \"20\u00d769=1380\"),\n    @(\"11\u00d790=990\", \"75\u00d779=5925\"),\n    @(\"44\u00d777=3388\", \"91\u00d796=8736\"),\n    @(\"47\u00d774=3478\", \"99\u00d773=7227\"),\n    @(\"65\u00d752=3380\", \"24\u00d760=1440\"),\n    @(\"46\u00d729=1334\", \"88\u00d755=4840\"),\n    @(\"90\u00d733=2970\", \"60\u00d746=2760\"),\n    @(\"46\u00d721=966\", \"90\u00d775=6750\"),\n    @(\"29\u00d720=580\", \"54\u00d733=1782\"),\n    @(\"56\u00d756=3136\", \"60\u00d788=5280\"),\n    @(\"77\u00d761=4697\", \"12\u00d783=996\"),\n    @(\"63\u00d737=2331\", \"48\u00d756=2688\"),\n    @(\"46\u00d727=1242\", \"71\u00d760=4260\"),\n    @(\"89\u00d714=1246\", \"42\u00d789=3738\"),\n    @(\"71\u00d787=6177\", \"62\u00d735=2170\"),\n    @(\"74\u00d740=2960\", \"99\u00d751=5049\"),\n    @(\"49\u00d751=2499\", \"32\u00d729=928\"),\n    @(\"45\u00d745=2025\", \"13\u00d794=1222\"),\n    @(\"62\u00d750=3100\", \"52\u00d711=572\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
